$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Longueur d'onde"
$ws.Range("B1").Value = "Absorbance de la solution diluée x100"

# Data rows: wavelength (A) vs absorbance (B)
$data = @(
    @(250, 0.143),
    @(260, 0.185),
    @(270, 0.275),
    @(280, 0.33),
    @(290, 0.21),
    @(300, 0.038),
    @(310, -0.006),
    @(320, -0.011),
    @(330, -0.01),
    @(340, -0.012),
    @(350, -0.012)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Select the last entered cell, matching the saved workbook's active selection
$ws.Range("B12").Select()
